$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.692.33'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '3.031.71'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''510.54'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").Value = '''140.23'
$ws.Range("E6").Value = '  +3.34%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.441'
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").Value = '''7.51'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  +3.06%  '
$ws.Range("E11").Value = '  +4.34%  '
$ws.Range("D12").Value = '3.559.26'
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '''26.65'
$ws.Range("E14").Value = '  +3.98%  '
$ws.Range("D15").Value = '''0.0000167'
$ws.Range("E15").Value = '  +8.58%  '
$ws.Range("D16").Value = '57.719.85'
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("D17").Value = '''6.25'
$ws.Range("E17").Value = '  +7.55%  '
$ws.Range("D18").Value = '3.031.14'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").Value = '''12.88'
$ws.Range("E19").Value = '  +3.71%  '
$ws.Range("D20").Value = '''8.02'
$ws.Range("E20").Value = '  +2.89%  '
$ws.Range("D21").Value = '''332.45'
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").Value = '''5.79'
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  +5.52%  '
$ws.Range("D25").Value = '''64.63'
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("D26").Value = '''0.169'
$ws.Range("E26").Value = '  +3.54%  '
$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  +3.99%  '
$ws.Range("D29").Value = '''6.81'
$ws.Range("E29").Value = '  +5.38%  '
$ws.Range("D30").Value = '''7.53'
$ws.Range("E30").Value = '  +9.62%  '
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("D32").Value = '''1.21'
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("D33").Value = '''20.75'
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("D34").Value = '''4.72'
$ws.Range("E34").Value = '  +5.45%  '
$ws.Range("D35").Value = '''155.25'
$ws.Range("E35").Value = '  -1.65%  '
$ws.Range("D36").Value = '''5.88'
$ws.Range("E36").Value = '  +6.08%  '
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").Value = '''24.91'
$ws.Range("E38").Value = '  +7.57%  '
$ws.Range("D39").Value = '''0.0686'
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").Value = '3.066.69'
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("D41").Value = '''37.50'
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").Value = '''3.88'
$ws.Range("E42").Value = '  +8.37%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '2.311.61'
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("E45").Value = '  +2.18%  '
$ws.Range("E46").Value = '  +2.15%  '
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("E48").Value = '  +4.73%  '
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("E50").Value = '  +3.72%  '
$ws.Range("E51").Value = '  -4.03%  '
